$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54 (shifts existing rows 54..143 down to 55..144,
# copying formatting from the row above as Excel normally does).
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new price record.
$ws.Range("A54").Value = 10
$ws.Range("B54").Value = "Vega Modelo de Temuco"
$ws.Range("C54").Value = "La Araucanía"
$ws.Range("D54").Value = "2021-09-30"
$ws.Range("E54").Value = 9
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100102
$ws.Range("H54").Value = "Cítricos"
$ws.Range("I54").Value = 100102006
$ws.Range("J54").Value = "Pomelo"
$ws.Range("K54").Value = "Start Ruby"
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 80
$ws.Range("N54").Value = 10000
$ws.Range("O54").Value = 10000
$ws.Range("P54").Value = 10000
$ws.Range("Q54").Value = "$/bandeja 15 kilos granel"
$ws.Range("R54").Value = "Región de O'Higgins"
$ws.Range("S54").Value = 667
$ws.Range("T54").Value = 15
